$wb = $excel.ActiveWorkbook

# --- Repayment Schedule sheet: insert a new (blank) column before the
#     existing "Late" column (column N), shifting Late/Outstanding right. ---
$ws = $wb.Worksheets.Item("Repayment Schedule")
$ws.Columns.Item(14).Insert()

# New column N is blank on every row -- nothing further to set there.

# --- Make "Repayment Schedule" the active sheet/tab, with cell R6 selected. ---
$ws.Activate()
$ws.Range("R6").Select()
